$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The Price (D) and Volume(1h) (E) columns are stored as plain Text
# in this sheet (prices use "." as a thousands separator, e.g.
# "29.447.71", which is not a valid number, and the percentages carry
# padding spaces) - every cell in both columns already has a Text
# format. Assigning a new value through Range.Value still lets Excel
# auto-detect numbers/dates from the string, which would corrupt
# values like "1.000" (-> 1) or "23.50" (-> 23.5) and flip the cell
# to the Number type. Pin every D/E cell we touch to Text first so
# the typed values round-trip exactly like the original sheet.
$textCells = "D2","E2","D3","E3","D4","E4","D5","E5","D6","E6","D7","D8","E8","D9","E9","D10","E11","D12","E12","D13","E13","D14","E14","D15","E15","D16","E16","D17","E17","D18","E18","D19","E19","D20","E20","E21","E22","D23","E23","E24","D25","E25","D26","E26","D27","E27","D28","E28","D29","E29","D30","E30","D31","E32","D33","E33","D34","E34","D35","E36","D37","E37","D38","E38","D39","E39","D40","E40","D41","E41","D42","E42","D43","E43","D44","E44","D45","E45","D46","E46","D47","E47","D48","E48","E49","D50","E50","D51","E51"
foreach ($addr in $textCells) { $ws.Range($addr).NumberFormat = "@" }

$ws.Range("D2").Value = "29.427.48"
$ws.Range("E2").Value = "  +1.93%  "
$ws.Range("D3").Value = "1.852.95"
$ws.Range("E3").Value = "  +1.14%  "
$ws.Range("D4").Value = "1.000"
$ws.Range("E4").Value = "  +0.05%  "
$ws.Range("D5").Value = "245.33"
$ws.Range("E5").Value = "  +0.03%  "
$ws.Range("D6").Value = "0.6928"
$ws.Range("E6").Value = "  +0.53%  "
$ws.Range("D7").Value = "1.001"
$ws.Range("D8").Value = "0.07667"
$ws.Range("E8").Value = "  -0.45%  "
$ws.Range("D9").Value = "0.3062"
$ws.Range("E9").Value = "  +0.46%  "
$ws.Range("D10").Value = "23.50"
$ws.Range("E11").Value = "  -0.62%  "
$ws.Range("B12").Value = "Polkadot"
$ws.Range("C12").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D12").Value = "5.146"
$ws.Range("E12").Value = "  +1.13%  "
$ws.Range("B13").Value = "WrappedEther"
$ws.Range("C13").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D13").Value = "1.847.95"
$ws.Range("E13").Value = "  +0.79%  "
$ws.Range("D14").Value = "0.6940"
$ws.Range("E14").Value = "  +1.93%  "
$ws.Range("D15").Value = "91.01"
$ws.Range("E15").Value = "  +0.14%  "
$ws.Range("D16").Value = "6.304"
$ws.Range("E16").Value = "  -2.06%  "
$ws.Range("D17").Value = "29.439.92"
$ws.Range("E17").Value = "  +2.00%  "
$ws.Range("D18").Value = "0.000008286"
$ws.Range("E18").Value = "  -0.32%  "
$ws.Range("D19").Value = "2.097.11"
$ws.Range("E19").Value = "  +1.20%  "
$ws.Range("D20").Value = "236.47"
$ws.Range("E20").Value = "  -2.29%  "
$ws.Range("E21").Value = "  +0.07%  "
$ws.Range("E22").Value = "  +0.07%  "
$ws.Range("D23").Value = "7.650"
$ws.Range("E23").Value = "  +2.63%  "
$ws.Range("E24").Value = "  +0.07%  "
$ws.Range("D25").Value = "0.1478"
$ws.Range("E25").Value = "  -0.13%  "
$ws.Range("D26").Value = "8.928"
$ws.Range("E26").Value = "  +1.56%  "
$ws.Range("D27").Value = "160.04"
$ws.Range("E27").Value = "  +1.04%  "
$ws.Range("D28").Value = "18.21"
$ws.Range("E28").Value = "  +0.03%  "
$ws.Range("D29").Value = "1.532"
$ws.Range("E29").Value = "  -0.88%  "
$ws.Range("D30").Value = "4.246"
$ws.Range("E30").Value = "  +0.69%  "
$ws.Range("D31").Value = "4.133"
$ws.Range("E32").Value = "  +0.97%  "
$ws.Range("D33").Value = "0.05204"
$ws.Range("E33").Value = "  +1.99%  "
$ws.Range("D34").Value = "0.7745"
$ws.Range("E34").Value = "  -0.24%  "
$ws.Range("D35").Value = "1.868"
$ws.Range("E36").Value = "  +0.48%  "
$ws.Range("D37").Value = "2.682"
$ws.Range("E37").Value = "  -0.31%  "
$ws.Range("D38").Value = "1.327.28"
$ws.Range("E38").Value = "  +8.92%  "
$ws.Range("D39").Value = "0.01868"
$ws.Range("E39").Value = "  +0.97%  "
$ws.Range("D40").Value = "2.722"
$ws.Range("E40").Value = "  +1.04%  "
$ws.Range("D41").Value = "0.9425"
$ws.Range("E41").Value = "  -1.75%  "
$ws.Range("D42").Value = "105.99"
$ws.Range("E42").Value = "  -2.40%  "
$ws.Range("D43").Value = "5.811"
$ws.Range("E43").Value = "  -0.33%  "
$ws.Range("D44").Value = "1.000"
$ws.Range("E44").Value = "  +0.12%  "
$ws.Range("D45").Value = "9.720"
$ws.Range("E45").Value = "  +0.93%  "
$ws.Range("B46").Value = "RocketPoolETH"
$ws.Range("C46").Value = "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
$ws.Range("D46").Value = "2.004.18"
$ws.Range("E46").Value = "  +1.46%  "
$ws.Range("B47").Value = "BabyDogeCoin"
$ws.Range("C47").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D47").Value = "0.00000000124"
$ws.Range("E47").Value = "  +1.19%  "
$ws.Range("D48").Value = "0.5224"
$ws.Range("E48").Value = "  +1.28%  "
$ws.Range("E49").Value = "  +2.13%  "
$ws.Range("D50").Value = "63.02"
$ws.Range("E50").Value = "  -2.06%  "
$ws.Range("D51").Value = "0.05956"
$ws.Range("E51").Value = "  +0.73%  "
